$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column S (19th column) -- shifts S..AN to T..AO
$ws.Columns("S").Insert()

# New column header
$ws.Range("S1").Value = "Sub brand"

# Restore autofilter over the new, wider range
$ws.AutoFilterMode = $false
$ws.Range("A1:AO35").AutoFilter() | Out-Null

# Update workbook-scoped defined names to the new range
foreach ($n in $wb.Names) {
    $n.RefersTo = "='Traditional Trade'!`$A`$1:`$AO`$35"
}

# Update the active selection to match the post-edit state
$ws.Range("S2").Select() | Out-Null
